# Weekly update: insert two new price records (date 44491) for
# "Agrícola del Norte S.A. de Arica" - Zapallo italiano, pushing the
# existing rows 192..218 down to 194..220.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at 192-193; everything from the old row 192
# onward shifts down by two rows (old 192 -> 194, ..., old 218 -> 220).
$ws.Range("A192:A193").EntireRow.Insert()

# New row 192: Primera quality entry for the week of 44491
$ws.Range("A192").Value = 1
$ws.Range("B192").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C192").Value = "Arica y Parinacota"
$ws.Range("D192").Value = 44491
$ws.Range("E192").Value = 15
$ws.Range("F192").Value = 100112032
$ws.Range("G192").Value = "Zapallo italiano"
$ws.Range("H192").Value = "Huracán"
$ws.Range("I192").Value = "Primera"
$ws.Range("J192").Value = 120
$ws.Range("K192").Value = 7000
$ws.Range("L192").Value = 8000
$ws.Range("M192").Value = 7500
$ws.Range("N192").Value = "`$/caja 70 unidades"
$ws.Range("O192").Value = "Región de Arica y Parinacota"
$ws.Range("P192").Value = 107
$ws.Range("Q192").Value = 70
$ws.Range("R192").Value = "Hortaliza"

# New row 193: Segunda quality entry for the week of 44491
$ws.Range("A193").Value = 1
$ws.Range("B193").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C193").Value = "Arica y Parinacota"
$ws.Range("D193").Value = 44491
$ws.Range("E193").Value = 15
$ws.Range("F193").Value = 100112032
$ws.Range("G193").Value = "Zapallo italiano"
$ws.Range("H193").Value = "Huracán"
$ws.Range("I193").Value = "Segunda"
$ws.Range("J193").Value = 160
$ws.Range("K193").Value = 5000
$ws.Range("L193").Value = 6000
$ws.Range("M193").Value = 5500
$ws.Range("N193").Value = "`$/caja 100 unidades"
$ws.Range("O193").Value = "Región de Arica y Parinacota"
$ws.Range("P193").Value = 55
$ws.Range("Q193").Value = 100
$ws.Range("R193").Value = "Hortaliza"
